$wb = $excel.ActiveWorkbook

# --- Sheet: "Generator Data" ---
$ws1 = $wb.Worksheets.Item("Generator Data")

$ws1.Range("B2").Value = 28062.4080285
$ws1.Range("C2").Value = 0.0553638950717

$ws1.Range("B3").Value = 5612.4816057
$ws1.Range("C3").Value = 0.06643667408604

$ws1.Range("B4").Value = 561.24816057
$ws1.Range("C4").Value = 0.029896503338718

$ws1.Range("B5").Value = 17824.5813342
$ws1.Range("C5").Value = 0.171523816216

# --- Sheet: "Yearly Fuel Costs" ---
$ws2 = $wb.Worksheets.Item("Yearly Fuel Costs")

$ws2.Range("B2").Value = 0.04125488801850326
$ws2.Range("B3").Value = 0.04125193280751012
$ws2.Range("B4").Value = 76.47893717906528
$ws2.Range("B5").Value = 81.92750187927129
$ws2.Range("B6").Value = 81.92750187927129
$ws2.Range("B7").Value = 81.92750187927129
$ws2.Range("B8").Value = 100.5919460951342
$ws2.Range("B9").Value = 3717.077487800779
$ws2.Range("B10").Value = 3739.768033320368
$ws2.Range("B11").Value = 3739.768033320368
$ws2.Range("B12").Value = 3739.768033320368
$ws2.Range("B13").Value = 3739.768033320368
$ws2.Range("B14").Value = 3739.768033320368
$ws2.Range("B15").Value = 3765.578516434562
$ws2.Range("B16").Value = 21976.74275005401
$ws2.Range("B17").Value = 22011.83316269384
$ws2.Range("B18").Value = 22011.83316269384
$ws2.Range("B19").Value = 22011.83316269384
$ws2.Range("B20").Value = 22011.83316269408
$ws2.Range("B21").Value = 22000.49476444961
